$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D category labels update (2022 -> 2023)
$catPib = "PIB 2023 Deflacionado"
$catVar = "Variação (%) 2023/2010"

# Row 2: DF
$ws.Range("B2").Value = 118174.1116095417
$ws.Range("D2").Value = $catPib

# Row 3: MT
$ws.Range("B3").Value = 76532.28963539573
$ws.Range("D3").Value = $catPib

# Row 4: SP
$ws.Range("B4").Value = 73845.19036585005
$ws.Range("D4").Value = $catPib

# Row 5: SC (was RJ)
$ws.Range("A5").Value = "SC"
$ws.Range("B5").Value = 69959.10090505773
$ws.Range("D5").Value = $catPib

# Row 6: RJ (was SC)
$ws.Range("A6").Value = "RJ"
$ws.Range("B6").Value = 67161.88535005336
$ws.Range("D6").Value = $catPib

# Row 7: MS
$ws.Range("B7").Value = 64948.89321994126
$ws.Range("D7").Value = $catPib

# Row 8: SE
$ws.Range("B8").Value = 26006.98661973922
$ws.Range("C8").Value = 23
$ws.Range("D8").Value = $catPib

# Row 9: BR
$ws.Range("B9").Value = 51300.70579350938
$ws.Range("D9").Value = $catPib

# Row 10: NE
$ws.Range("B10").Value = 26237.41536180414
$ws.Range("D10").Value = $catPib

# Row 11: PI (was MT)
$ws.Range("A11").Value = "PI"
$ws.Range("B11").Value = 1.537570136346218
$ws.Range("D11").Value = $catVar

# Row 12: AL (was PI)
$ws.Range("A12").Value = "AL"
$ws.Range("B12").Value = 1.47614098883597
$ws.Range("D12").Value = $catVar

# Row 13: BA
$ws.Range("B13").Value = 1.456636142415122
$ws.Range("D13").Value = $catVar

# Row 14: PR (was MS)
$ws.Range("A14").Value = "PR"
$ws.Range("B14").Value = 1.361548723808843
$ws.Range("D14").Value = $catVar

# Row 15: MT (was PR)
$ws.Range("A15").Value = "MT"
$ws.Range("B15").Value = 1.347864888435147
$ws.Range("D15").Value = $catVar

# Row 16: RS (was RJ)
$ws.Range("A16").Value = "RS"
$ws.Range("B16").Value = 1.330588661237085
$ws.Range("D16").Value = $catVar

# Row 17: SE
$ws.Range("B17").Value = 0.91785810919599
$ws.Range("C17").Value = 23
$ws.Range("D17").Value = $catVar

# Row 18: BR
$ws.Range("B18").Value = 1.205735709293767
$ws.Range("D18").Value = $catVar

# Row 19: NE
$ws.Range("B19").Value = 1.297733510014661
$ws.Range("D19").Value = $catVar
